# Finish schematics BOM: fill in Fabricant / Fournisseur / RefFabricant / RefFournisseur
# for the two Green/Red LED rows (D1,D3 and D2) that were previously left blank.
#
# A leading apostrophe forces Excel to store the value as text (so purely-numeric
# reference codes like "2846598" are kept as shared-string text, not numbers) while
# preserving the existing cell style (quote-prefix style already applied to these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: D1, D3 - Green LED (KINGBRIGHT / LED-0 / Farnell / KPTD-2012LVZGCK / 2846598)
$ws.Range("D3").Value = "'KINGBRIGHT"
$ws.Range("F3").Value = "'Farnell"
$ws.Range("I3").Value = "'KPTD-2012LVZGCK"
$ws.Range("J3").Value = "'2846598"

# Row 4: D2 - Red LED (KINGBRIGHT / LED-0 / Farnell / KPTD-2012LVSURCK / 2846595)
$ws.Range("D4").Value = "'KINGBRIGHT"
$ws.Range("F4").Value = "'Farnell"
$ws.Range("I4").Value = "'KPTD-2012LVSURCK"
$ws.Range("J4").Value = "'2846595"
